# Change the table style ("Table Style" gallery selection) applied to the
# table on slide 16 from its original style to the new style.
#
# The table lives in the graphic frame that is the 3rd shape on slide 16
# (shape 1 = title placeholder, shape 2 = picture, shape 3 = the table).
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)

if ($sh.HasTable) {
    $tbl = $sh.Table
    # Table styles are GUID-identified entries from the built-in Table
    # Styles gallery; they must be applied via ApplyStyle (a direct
    # assignment to .Style is rejected by PowerPoint).
    $tbl.ApplyStyle("{3286A7B5-318E-4AB5-8A84-EBAC4C7A7C66}")
}
